$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2069970845481049
$ws.Range("C2").Value = 0.5510204081632653
$ws.Range("J2").Value = 0.008746355685131196
$ws.Range("P2").Value = 0.1341107871720117
$ws.Range("S2").Value = 0.09912536443148688
$ws.Range("C3").Value = 0.005235602094240838
$ws.Range("J3").Value = 0.07329842931937172
$ws.Range("P3").Value = 0.6701570680628273
$ws.Range("S3").Value = 0.2513089005235602
$ws.Range("J4").Value = 0.09803921568627451
$ws.Range("P4").Value = 0.6078431372549019
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.04602510460251046
$ws.Range("D6").Value = 0.01255230125523013
$ws.Range("F6").Value = 0.04602510460251046
$ws.Range("J6").Value = 0.2887029288702929
$ws.Range("O6").Value = 0.02928870292887029
$ws.Range("Q6").Value = 0.1255230125523012
$ws.Range("R6").Value = 0.07949790794979079
$ws.Range("S6").Value = 0.3723849372384937
$ws.Range("B7").Value = 0.1443850267379679
$ws.Range("D7").Value = 0.0053475935828877
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1711229946524064
$ws.Range("O7").Value = 0.0213903743315508
$ws.Range("Q7").Value = 0.1764705882352941
$ws.Range("R7").Value = 0.06417112299465241
$ws.Range("S7").Value = 0.3582887700534759
$ws.Range("B8").Value = 0.1174004192872117
$ws.Range("D8").Value = 0.01886792452830189
$ws.Range("F8").Value = 0.06498951781970649
$ws.Range("J8").Value = 0.129979035639413
$ws.Range("O8").Value = 0.0209643605870021
$ws.Range("Q8").Value = 0.1740041928721174
$ws.Range("R8").Value = 0.08595387840670859
$ws.Range("S8").Value = 0.3878406708595388
$ws.Range("B9").Value = 0.1221374045801527
$ws.Range("D9").Value = 0.01526717557251908
$ws.Range("F9").Value = 0.0916030534351145
$ws.Range("J9").Value = 0.0916030534351145
$ws.Range("O9").Value = 0.01526717557251908
$ws.Range("Q9").Value = 0.1679389312977099
$ws.Range("R9").Value = 0.1068702290076336
$ws.Range("S9").Value = 0.3893129770992366
$ws.Range("B10").Value = 0.119969040247678
$ws.Range("D10").Value = 0.02863777089783282
$ws.Range("E10").Value = 0.003869969040247678
$ws.Range("F10").Value = 0.0781733746130031
$ws.Range("J10").Value = 0.1369969040247678
$ws.Range("O10").Value = 0.01934984520123839
$ws.Range("Q10").Value = 0.1818885448916409
$ws.Range("R10").Value = 0.07739938080495357
$ws.Range("S10").Value = 0.3537151702786377
$ws.Range("G11").Value = 0.1516245487364621
$ws.Range("J11").Value = 0.1046931407942238
$ws.Range("K11").Value = 0.2021660649819494
$ws.Range("L11").Value = 0.5234657039711191
$ws.Range("S11").Value = 0.01805054151624549
$ws.Range("G12").Value = 0.7712418300653595
$ws.Range("J12").Value = 0.1633986928104575
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.0457516339869281
$ws.Range("S12").Value = 0.0130718954248366
$ws.Range("G13").Value = 0.6730769230769231
$ws.Range("J13").Value = 0.2692307692307692
$ws.Range("S13").Value = 0.0576923076923077
$ws.Range("F15").Value = 0.0179372197309417
$ws.Range("H15").Value = 0.2062780269058296
$ws.Range("I15").Value = 0.08520179372197309
$ws.Range("J15").Value = 0.3318385650224215
$ws.Range("K15").Value = 0.05381165919282511
$ws.Range("M15").Value = 0.004484304932735426
$ws.Range("N15").Value = 0.008968609865470852
$ws.Range("O15").Value = 0.04484304932735426
$ws.Range("S15").Value = 0.2466367713004484
$ws.Range("F16").Value = 0.02955665024630542
$ws.Range("H16").Value = 0.1921182266009852
$ws.Range("I16").Value = 0.06403940886699508
$ws.Range("J16").Value = 0.4187192118226601
$ws.Range("K16").Value = 0.09852216748768473
$ws.Range("M16").Value = 0.02955665024630542
$ws.Range("O16").Value = 0.08374384236453201
$ws.Range("S16").Value = 0.08374384236453201
$ws.Range("F17").Value = 0.015
$ws.Range("H17").Value = 0.22
$ws.Range("I17").Value = 0.0525
$ws.Range("J17").Value = 0.4075
$ws.Range("K17").Value = 0.1325
$ws.Range("M17").Value = 0.02
$ws.Range("O17").Value = 0.065
$ws.Range("S17").Value = 0.08749999999999999
$ws.Range("F18").Value = 0.02197802197802198
$ws.Range("H18").Value = 0.1978021978021978
$ws.Range("I18").Value = 0.06593406593406594
$ws.Range("J18").Value = 0.467032967032967
$ws.Range("K18").Value = 0.0989010989010989
$ws.Range("M18").Value = 0.02747252747252747
$ws.Range("O18").Value = 0.07142857142857142
$ws.Range("S18").Value = 0.04945054945054945
$ws.Range("F19").Value = 0.0226890756302521
$ws.Range("H19").Value = 0.2260504201680672
$ws.Range("I19").Value = 0.05630252100840336
$ws.Range("J19").Value = 0.3857142857142857
$ws.Range("K19").Value = 0.09327731092436975
$ws.Range("M19").Value = 0.02941176470588235
$ws.Range("O19").Value = 0.073109243697479
$ws.Range("S19").Value = 0.1134453781512605
